{"js": "// Replace the date line and each two-digit-by-two-digit multiplication\n// prompt in the practice-sheet table with the new values from the\n// target revision. Every \"before\" string in the sheet is unique, so a\n// plain text search-and-replace (matching the whole cell text,\n// including the trailing \"=\") unambiguously targets the correct run.\nconst replacements = [\n  [\"2025-02-16 Sunday\", \"2025-02-17 Monday\"],\n  [\"88\u00d750=\", \"16\u00d778=\"],\n  [\"25\u00d752=\", \"23\u00d793=\"],\n  [\"96\u00d760=\", \"50\u00d738=\"],\n  [\"28\u00d790=\", \"14\u00d726=\"],\n  [\"92\u00d787=\", \"40\u00d762=\"],\n  [\"94\u00d722=\", \"68\u00d799=\"],\n  [\"44\u00d713=\", \"98\u00d726=\"],\n  [\"73\u00d733=\", \"86\u00d762=\"],\n  [\"12\u00d786=\", \"56\u00d763=\"],\n  [\"58\u00d778=\", \"40\u00d776=\"],\n  [\"57\u00d781=\", \"69\u00d742=\"],\n  [\"42\u00d778=\", \"11\u00d765=\"],\n  [\"84\u00d799=\", \"60\u00d754=\"],\n  [\"41\u00d794=\", \"95\u00d769=\"],\n  [\"32\u00d745=\", \"15\u00d797=\"],\n  [\"81\u00d759=\", \"39\u00d717=\"],\n  [\"69\u00d750=\", \"42\u00d777=\"],\n  [\"48\u00d776=\", \"18\u00d798=\"],\n  [\"70\u00d711=\", \"26\u00d795=\"],\n  [\"57\u00d770=\", \"90\u00d772=\"],\n  [\"56\u00d784=\", \"52\u00d712=\"],\n  [\"52\u00d746=\", \"63\u00d753=\"],\n  [\"86\u00d721=\", \"35\u00d762=\"],\n  [\"67\u00d736=\", \"92\u00d739=\"],\n  [\"27\u00d731=\", \"64\u00d722=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the practice-sheet date line and every two-digit-by-two-digit\n# multiplication prompt to the new values from the target revision.\n# Every \"before\" cell value in this sheet is unique, so a plain\n# Find/Replace on the whole-cell text (including the trailing \"=\")\n# unambiguously targets the correct run without disturbing any other\n# formatting (font, size, paragraph alignment, table structure, etc.).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-02-16 Sunday\", \"2025-02-17 Monday\"),\n  @(\"88\u00d750=\", \"16\u00d778=\"),\n  @(\"25\u00d752=\", \"23\u00d793=\"),\n  @(\"96\u00d760=\", \"50\u00d738=\"),\n  @(\"28\u00d790=\", \"14\u00d726=\"),\n  @(\"92\u00d787=\", \"40\u00d762=\"),\n  @(\"94\u00d722=\", \"68\u00d799=\"),\n  @(\"44\u00d713=\", \"98\u00d726=\"),\n  @(\"73\u00d733=\", \"86\u00d762=\"),\n  @(\"12\u00d786=\", \"56\u00d763=\"),\n  @(\"58\u00d778=\", \"40\u00d776=\"),\n  @(\"57\u00d781=\", \"69\u00d742=\"),\n  @(\"42\u00d778=\", \"11\u00d765=\"),\n  @(\"84\u00d799=\", \"60\u00d754=\"),\n  @(\"41\u00d794=\", \"95\u00d769=\"),\n  @(\"32\u00d745=\", \"15\u00d797=\"),\n  @(\"81\u00d759=\", \"39\u00d717=\"),\n  @(\"69\u00d750=\", \"42\u00d777=\"),\n  @(\"48\u00d776=\", \"18\u00d798=\"),\n  @(\"70\u00d711=\", \"26\u00d795=\"),\n  @(\"57\u00d770=\", \"90\u00d772=\"),\n  @(\"56\u00d784=\", \"52\u00d712=\"),\n  @(\"52\u00d746=\", \"63\u00d753=\"),\n  @(\"86\u00d721=\", \"35\u00d762=\"),\n  @(\"67\u00d736=\", \"92\u00d739=\"),\n  @(\"27\u00d731=\", \"64\u00d722=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
